$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold, border, centered) from H1 to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 and IF numeric data for rows 2-19
$data = @(
    @(9,9),
    @(4,5),
    @(9,9),
    @(7,7),
    @(8,8),
    @(7,8),
    @(9,9),
    @(8,8),
    @(8,9),
    @(9,9),
    @(8,8),
    @(5,6),
    @(8,9),
    @(5,6),
    @(8,8),
    @(6,6),
    @(5,5),
    @(7,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
